# Refresh the "Chart" sheet's rolling date window for the latest GSC export:
# drop the oldest date row and append the next day, so every existing row's
# Date/Items slide up by one and a new trailing row is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastDataRow = 90

# Remember the last row's current date + items count before the shift so we
# can compute the new trailing day and carry its count forward.
$lastDateText = $ws.Cells.Item($lastDataRow, 1).Value2
$lastItems = $ws.Cells.Item($lastDataRow, 3).Value2

# Drop the oldest day: deleting row 2 shifts rows 3..90 up into 2..89,
# carrying each cell's existing value/type with it (no retyping needed).
$ws.Rows.Item(2).Delete()

# Append the new trailing day at row 90. Pre-format as text so the
# ISO-formatted date string isn't auto-converted into a date serial.
$lastDate = [DateTime]::ParseExact($lastDateText, "yyyy-MM-dd", $null)
$nextDateText = $lastDate.AddDays(1).ToString("yyyy-MM-dd")

$ws.Cells.Item($lastDataRow, 1).NumberFormat = "@"
$ws.Cells.Item($lastDataRow, 1).Value = $nextDateText
$ws.Cells.Item($lastDataRow, 2).Value = 0.0
$ws.Cells.Item($lastDataRow, 3).Value = $lastItems
